$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the instruction / validation messages in rows 1-3
$ws.Range("A1").Value = "Kecuali kolom Gambar, semua kolom wajib diisi!"
$ws.Range("A3").Value = "Kosongi kolom Gambar, apabila soal terdapat Gambar bisa ditambahkan saat edit soal"

# Remove the two sample data rows (rows 5 and 6), leaving just the header row (row 4)
$ws.Range("A5:K6").Delete()

# Update the active selection to match the target workbook
$ws.Range("C8").Select()
